$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E1").Value = "La Liga"
$ws.Range("F1").Value = "Premier League"

$ws.Range("F2").Select()
